$wb = $excel.ActiveWorkbook

# ---- Build the first sheet (ArcFace) in place: rename + data + formatting ----
$wsArc = $wb.Worksheets.Item(1)
$wsArc.Name = "ArcFace"

$wsArc.Range("A1").Value = "Metric"
$wsArc.Range("B1").Value = "Value (Weighted)"
$wsArc.Range("C1").Value = "Value (Micro)"
$wsArc.Range("D1").Value = "Value (Macro)"

$wsArc.Range("A2").Value = "Accuracy"
$wsArc.Range("B2").Value = 0.99673400000000001
$wsArc.Range("C2").Value = 0.99673400000000001
$wsArc.Range("D2").Value = 0.99673400000000001

$wsArc.Range("A3").Value = "Precision"
$wsArc.Range("B3").Value = 0.99790699999999999
$wsArc.Range("C3").Value = 0.99673400000000001
$wsArc.Range("D3").Value = 0.79874500000000004

$wsArc.Range("A4").Value = "Recall"
$wsArc.Range("B4").Value = 0.99673400000000001
$wsArc.Range("C4").Value = 0.99673400000000001
$wsArc.Range("D4").Value = 0.99745200000000001

$wsArc.Range("A5").Value = "F1-Score"
$wsArc.Range("B5").Value = 0.99715200000000004
$wsArc.Range("C5").Value = 0.99673400000000001
$wsArc.Range("D5").Value = 0.85524

# column B:D sized to fit their header/content (column A keeps its original width)
$wsArc.Range("B1:D5").EntireColumn.AutoFit()
$wsArc.Columns.Item(2).ColumnWidth = 15.26953125
$wsArc.Columns.Item(3).ColumnWidth = 12
$wsArc.Columns.Item(4).ColumnWidth = 12.54296875

$wsArc.PageSetup.PaperSize = 9
$wsArc.PageSetup.Orientation = 1

# ---- Duplicate ArcFace -> VGGFace -> FaceNet512 so formatting/column widths carry over exactly ----
$wsArc.Copy($null, $wsArc)
$wsVgg = $wb.Worksheets.Item(2)
$wsVgg.Name = "VGGFace"

$wsVgg.Copy($null, $wsVgg)
$wsFace = $wb.Worksheets.Item(3)
$wsFace.Name = "FaceNet512"

# ---- VGGFace + FaceNet512 share the same metric values, different from ArcFace ----
foreach ($ws in @($wsVgg, $wsFace)) {
    $ws.Range("A2").Value = "Accuracy"
    $ws.Range("B2").Value = 0.99738700000000002
    $ws.Range("C2").Value = 0.99738700000000002
    $ws.Range("D2").Value = 0.99738700000000002

    $ws.Range("A3").Value = "Precision"
    $ws.Range("B3").Value = 0.99912900000000004
    $ws.Range("C3").Value = 0.99738700000000002
    $ws.Range("D3").Value = 0.77777799999999997

    $ws.Range("A4").Value = "Recall"
    $ws.Range("B4").Value = 0.99738700000000002
    $ws.Range("C4").Value = 0.99738700000000002
    $ws.Range("D4").Value = 0.99749399999999999

    $ws.Range("A5").Value = "F1-Score"
    $ws.Range("B5").Value = 0.99803600000000003
    $ws.Range("C5").Value = 0.99738700000000002
    $ws.Range("D5").Value = 0.83207500000000001
}

# ---- Per-sheet selection state (persists in the saved file even for inactive sheets) ----
$wsArc.Range("E1").Select()
$wsVgg.Range("E3").Select()
$wsFace.Range("B2").Select()
